$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.058.67'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '1.619.15'
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.45'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.06'
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").Value = '1.624.11'
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.15'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.540'
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.66'
$ws.Range("E15").Value = '  -4.08%  '
$ws.Range("D16").Value = '27.034.90'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.46'
$ws.Range("E18").Value = '  -3.06%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.87'
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.40'
$ws.Range("E22").Value = '  -5.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.07'
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.29'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.43'
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.53'
$ws.Range("E28").Value = '  -1.76%  '
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.759'
$ws.Range("E31").Value = '  +37.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.37'
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.01'
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").Value = '1.343.34'
$ws.Range("E34").Value = '  +3.56%  '
$ws.Range("E35").Value = '  -1.51%  '
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("E37").Value = '  +0.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.853'
$ws.Range("E38").Value = '  -1.94%  '
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.800'
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '65.26'
$ws.Range("E41").Value = '  +4.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.22'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("D44").Value = '1.755.84'
$ws.Range("E44").Value = '  -1.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.899'
$ws.Range("E45").Value = '  +33.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.03'
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0516'
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.65'
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.06%  '
